# Kiviat_Template.xlsx — "Add files via upload" edit
#
# Functional changes applied through the Excel object model:
#   1. Update the Score Page A / Score Page B figures for several of the
#      Golden Rules rows (B4:D11 table that feeds the radar/kiviat chart).
#   2. Turn the two URL cells (C2/D2) into real hyperlinks (they already
#      held the URL text; Excel re-styles them with the built-in
#      "Hyperlink" cell style once a hyperlink is attached).
#   3. Move the active selection to B11 (last table row) instead of D10,
#      matching the cursor position the author left the sheet in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Updated scores -----------------------------------------------
# row 5  (2. Seek Universal Usability): Score Page A 9 -> 8
$ws.Range("C5").Value = 8
# row 6  (3. Offer Informative Fedback): Score Page A 8 -> 6
$ws.Range("C6").Value = 6
# row 8  (5. Prevent Errors): Score Page A 8 -> 10, Score Page B 4 -> 10
$ws.Range("C8").Value = 10
$ws.Range("D8").Value = 10
# row 10 (7. Keep Users in Control): Score Page A 10 -> 2, Score Page B 10 -> 0
$ws.Range("C10").Value = 2
$ws.Range("D10").Value = 0

# --- 2. Hyperlinks on the two source-URL cells ------------------------
$ws.Hyperlinks.Add($ws.Range("C2"), "https://cabildo.grancanaria.com/")
$ws.Hyperlinks.Add($ws.Range("D2"), "https://www.tenerife.es/portalcabtfe/es/")

# --- 3. Selection / cursor position ------------------------------------
$ws.Range("B11").Select() | Out-Null
